# Bump the document revision marker from "...2020.1-EN-rev1" to
# "...2020.1-EN-rev2" (jira-FL-33 / 20-04-36 distractor fix, rev2).
#
# The digit lives in its own run immediately after a run ending in
# "...2020.1-EN-rev", both on the title page (main document story) and
# in the page footer. We locate the anchor text with Find, then replace
# only the single trailing digit so the surrounding runs/formatting are
# left untouched.

$d = $word.ActiveDocument

function Bump-RevisionDigit($story) {
    $anchor = $story.Duplicate
    $found = $anchor.Find.Execute("2020.1-EN-rev", $true, $false, $false, `
                                   $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }

    $digit = $story.Duplicate
    $digit.SetRange($anchor.End, $anchor.End + 1)
    if ($digit.Text -eq "1") {
        $digit.Text = "2"
        return $true
    }
    return $false
}

# Main document body (title-page "Document version: 2020.1-EN-rev1,").
Bump-RevisionDigit($d.Content) | Out-Null

# Default page footer ("Version 2020.1-EN-rev1").
$section = $d.Sections(1)
$footer = $section.Footers(1)
Bump-RevisionDigit($footer.Range) | Out-Null
